$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "United Arab Emirates"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "4"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = "Insurance (Prop/Cas.)"
$ws.Range("D2").Value = 0.08595
$ws.Range("E2").Value = 0.0984
$ws.Range("G2").Value = 0.07426844014510278
$ws.Range("H2").Value = 0.07426844014510278
$ws.Range("I2").Value = 0.1384662334291068
$ws.Range("J2").Value = 0.1355354617866735
$ws.Range("K2").Value = 53.45
$ws.Range("L2").Value = 0.1292623941958888
$ws.Range("M2").Value = 11.555
$ws.Range("N2").Value = 0.02880109670987039
$ws.Range("O2").Value = 0.2161833489242283
$ws.Range("P2").Value = 11.302
$ws.Range("Q2").Value = 0.02817048853439681
$ws.Range("R2").Value = 0.2114499532273152
$ws.Range("S2").Value = 0.253
$ws.Range("T2").Value = 0.02189528342708784
$ws.Range("U2").Value = 118.1
$ws.Range("V2").Value = 0.2943668993020937
$ws.Range("W2").Value = 0.1694444444444445
$ws.Range("X2").Value = 0.044565396872846
$ws.Range("Y2").Value = 0.1248790475715985
$ws.Range("Z2").Value = 2.404334490465921
$ws.Range("AA2").Value = 0.1420746779136872
$ws.Range("AB2").Value = 0.04246886801701201
$ws.Range("AC2").Value = 0.1032551068660386
$ws.Range("AD2").Value = 73.68000000000001
$ws.Range("AE2").Value = 0.03106238532161488
$ws.Range("AF2").Value = 73.71106238532161
$ws.Range("AG2").Value = -44.38893761467838
$ws.Range("AH2").Value = 0.1552102450827219
$ws.Range("AI2").Value = 0.1744182037480943
$ws.Range("AJ2").Value = -0.1244046003448811
$ws.Range("AK2").Value = -0.1457711823897865
$ws.Range("AL2").Value = 0.949
$ws.Range("AM2").Value = 0.949
$ws.Range("AN2").Value = 1.215941909398465
$ws.Range("AO2").Value = 60.30558482613277
$ws.Range("AP2").Value = -0.7325511612291177
$ws.Range("AQ2").Value = 60.30558482613277

# Row 3
$ws.Range("A3").Value = "United Arab Emirates"
$ws.Range("B3").Value = "Insurance House P.S.C. (ADX:IH)"
$ws.Range("C3").Value = "Insurance (Prop/Cas.)"
$ws.Range("D3").Value = 0.134
$ws.Range("G3").Value = 0.1204641350210971
$ws.Range("H3").Value = 0.1204641350210971
$ws.Range("I3").Value = 0.09156118143459915
$ws.Range("J3").Value = 0.09156118143459915
$ws.Range("K3").Value = 3.69
$ws.Range("L3").Value = 0.07784810126582278
$ws.Range("M3").Value = 0.515
$ws.Range("N3").Value = 0.01973180076628352
$ws.Range("O3").Value = 0.1395663956639566
$ws.Range("P3").Value = 0.262
$ws.Range("Q3").Value = 0.01003831417624521
$ws.Range("R3").Value = 0.07100271002710028
$ws.Range("S3").Value = 0.253
$ws.Range("T3").Value = 0.4912621359223301
$ws.Range("U3").Value = 10.5
$ws.Range("V3").Value = 0.4022988505747126
$ws.Range("W3").Value = 0.1078947368421053
$ws.Range("X3").Value = 0.04334820328349445
$ws.Range("Y3").Value = 0.0645465335586108
$ws.Range("Z3").Value = 4.232142857142856
$ws.Range("AA3").Value = 0.3874999999999999
$ws.Range("AB3").Value = 0.04334820328349445
$ws.Range("AC3").Value = 0.3441517967165054
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -10.5
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.673076923076923
$ws.Range("AK3").Value = -0.4285714285714285
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AN3").Value = 0
$ws.Range("AP3").Value = -2.164948453608248

# Row 4
$ws.Range("A4").Value = "United Arab Emirates"
$ws.Range("B4").Value = "Islamic Arab Insurance Co. (Salama) PJSC (DFM:SALAMA)"
$ws.Range("C4").Value = "Insurance (Prop/Cas.)"
$ws.Range("D4").Value = 0.0379
$ws.Range("G4").Value = 0.07788944723618091
$ws.Range("H4").Value = 0.07788944723618091
$ws.Range("I4").Value = 0.2153508690240187
$ws.Range("J4").Value = 0.1971184320037079
$ws.Range("K4").Value = 40.9
$ws.Range("L4").Value = 0.1712730318257956
$ws.Range("M4").Value = 9.699999999999999
$ws.Range("N4").Value = 0.03544026306174643
$ws.Range("O4").Value = 0.2371638141809291
$ws.Range("P4").Value = 9.699999999999999
$ws.Range("Q4").Value = 0.03544026306174643
$ws.Range("R4").Value = 0.2371638141809291
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 61.6
$ws.Range("V4").Value = 0.2250639386189258
$ws.Range("W4").Value = 0.2010816125860374
$ws.Range("X4").Value = 0.04346134210090338
$ws.Range("Y4").Value = 0.157620270485134
$ws.Range("Z4").Value = 1.788347937432797
$ws.Range("AA4").Value = 0.352516341303818
$ws.Range("AB4").Value = 0.04324233239675908
$ws.Range("AC4").Value = 0.3092740089070589
$ws.Range("AD4").Value = 1.99
$ws.Range("AE4").Value = 0.03106238532161488
$ws.Range("AF4").Value = 2.021062385321615
$ws.Range("AG4").Value = -59.57893761467839
$ws.Range("AH4").Value = 0.007330097917935517
$ws.Range("AI4").Value = 0.007884886113195792
$ws.Range("AJ4").Value = -0.2782488418045633
$ws.Range("AK4").Value = -0.3059706889683113
$ws.Range("AL4").Value = 0.949
$ws.Range("AM4").Value = 0.949
$ws.Range("AN4").Value = 0.03752451350128225
$ws.Range("AO4").Value = 54.1622760800843
$ws.Range("AP4").Value = -1.123452587393996
$ws.Range("AQ4").Value = 54.1622760800843
$ws.Range("E4").ClearContents()

# Row 5
$ws.Range("A5").Value = "United Arab Emirates"
$ws.Range("B5").Value = "Dar Al Takaful PJSC (DFM:DARTAKAFUL)"
$ws.Range("C5").Value = "Insurance (Prop/Cas.)"
$ws.Range("D5").Value = 0.216
$ws.Range("G5").Value = 0.03561497326203208
$ws.Range("H5").Value = 0.03561497326203208
$ws.Range("I5").Value = -0.02064171122994652
$ws.Range("J5").Value = -0.02064171122994652
$ws.Range("K5").Value = 6.13
$ws.Range("L5").Value = 0.06556149732620321
$ws.Range("M5").Value = 1.34
$ws.Range("N5").Value = 0.03517060367454068
$ws.Range("O5").Value = 0.2185970636215335
$ws.Range("P5").Value = 1.34
$ws.Range("Q5").Value = 0.03517060367454068
$ws.Range("R5").Value = 0.2185970636215335
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 18.2
$ws.Range("V5").Value = 0.4776902887139107
$ws.Range("W5").Value = 0.1792397660818713
$ws.Range("X5").Value = 0.06832135395953781
$ws.Range("Y5").Value = 0.1109184121223335
$ws.Range("Z5").Value = 3.312079348211122
$ws.Range("AA5").Value = -0.06836698547644349
$ws.Range("AB5").Value = 0.03439680969853831
$ws.Range("AC5").Value = -0.1027637951749818
$ws.Range("AD5").Value = 62.1
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 62.1
$ws.Range("AG5").Value = 43.90000000000001
$ws.Range("AH5").Value = 0.6197604790419161
$ws.Range("AI5").Value = 0.6142433234421365
$ws.Range("AJ5").Value = 0.5353658536585366
$ws.Range("AK5").Value = 0.5295536791314838
$ws.Range("AL5").Value = 0
$ws.Range("AM5").Value = 0
$ws.Range("AN5").Value = -82.03434610303832
$ws.Range("AP5").Value = -57.99207397622194
$ws.Range("E5").ClearContents()
$ws.Range("AO5").ClearContents()
$ws.Range("AQ5").ClearContents()

# Row 6
$ws.Range("A6").Value = "United Arab Emirates"
$ws.Range("B6").Value = "Dubai Islamic Insurance & Reinsurance Co. (Aman) (P.J.S.C) (DFM:AMAN)"
$ws.Range("C6").Value = "Insurance (Prop/Cas.)"
$ws.Range("D6").Value = -0.09119999999999999
$ws.Range("E6").Value = 0.0984
$ws.Range("G6").Value = 0.0908284023668639
$ws.Range("H6").Value = 0.0908284023668639
$ws.Range("I6").Value = 0.1011834319526627
$ws.Range("J6").Value = 0.1011834319526627
$ws.Range("K6").Value = 2.73
$ws.Range("L6").Value = 0.08076923076923077
$ws.Range("M6").Value = -0
$ws.Range("N6").Value = -0
$ws.Range("O6").Value = -0
$ws.Range("P6").Value = -0
$ws.Range("Q6").Value = -0
$ws.Range("R6").Value = -0
$ws.Range("S6").Value = 0
$ws.Range("U6").Value = 27.8
$ws.Range("V6").Value = 0.439178515007899
$ws.Range("W6").Value = 0.1596491228070175
$ws.Range("X6").Value = 0.04566945164478861
$ws.Range("Y6").Value = 0.1139796711622289
$ws.Range("Z6").Value = -34.48979591836733
$ws.Range("AA6").Value = -3.489795918367345
$ws.Range("AB6").Value = 0.04169540363726493
$ws.Range("AC6").Value = -3.53149132200461
$ws.Range("AD6").Value = 9.59
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 9.59
$ws.Range("AG6").Value = -18.21
$ws.Range("AH6").Value = 0.1315681163396899
$ws.Range("AI6").Value = 0.3176548526001987
$ws.Range("AJ6").Value = -0.4038589487691285
$ws.Range("AK6").Value = -7.619246861924685
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AN6").Value = 2.763688760806916
$ws.Range("AP6").Value = -5.247838616714698

Write-Host "Edit applied"